# Drugs@FDA reference sheet: remove the
# submissions.submission_property_type.id row (#207) - it's an internal
# primary key without a public meaning.
#
# In the "device_classification_fields" sheet, row 43 holds:
#   A: submission_property_type
#   B: id
#   C: string
#   D: The id of the submission property type.
#
# Deleting that whole row shifts every subsequent row up by one (so the
# used range shrinks from A1:D47 to A1:D46) and drops the now-unused
# "id" / "The id of the submission property type." shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Rows(43)

# Mirror the interactive edit: select the row before removing it, so the
# saved selection reflects the cell where the deletion happened.
[void]$row.EntireRow.Select()
[void]$row.EntireRow.Delete()

# The session was also zoomed in to 130% when the edit was made.
$excel.ActiveWindow.Zoom = 130
